$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price strings so they are not
# auto-coerced into numbers by Excel (matches source data being stored as text).
$textCells = @("D5", "D7", "D8", "D9", "D10", "D11", "D14", "D15", "D16", "D18", "D21", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D32", "D33", "D34", "D35", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data (prices / volume %) to reflect the latest scrape.
$ws.Range("D2").Value = "30.775.30"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "1.928.17"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("D5").Value = "242.37"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "0.4858"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").Value = "0.2933"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("D9").Value = "0.06835"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").Value = "19.17"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").Value = "106.45"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("D13").Value = "1.928.00"
$ws.Range("E13").Value = "  -1.36%  "
$ws.Range("D14").Value = "5.325"
$ws.Range("E14").Value = "  -2.35%  "
$ws.Range("D15").Value = "0.6961"
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").Value = "275.11"
$ws.Range("E16").Value = "  -2.65%  "
$ws.Range("D17").Value = "30.753.78"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "0.000007667"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("D21").Value = "12.97"
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "6.449"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("D24").Value = "9.851"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").Value = "164.97"
$ws.Range("E25").Value = "  -2.97%  "
$ws.Range("D26").Value = "19.44"
$ws.Range("E26").Value = "  -2.58%  "
$ws.Range("D27").Value = "2.147"
$ws.Range("E27").Value = "  -2.20%  "
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("D29").Value = "1.382"
$ws.Range("E29").Value = "  -1.97%  "
$ws.Range("D30").Value = "4.565"
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("E31").Value = "  -2.08%  "
$ws.Range("D32").Value = "4.353"
$ws.Range("E32").Value = "  -2.02%  "
$ws.Range("D33").Value = "0.04874"
$ws.Range("E33").Value = "  -1.10%  "
$ws.Range("D34").Value = "0.7571"
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("D35").Value = "1.140"
$ws.Range("E35").Value = "  -2.51%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "2.723"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("D39").Value = "2.643"
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "77.48"
$ws.Range("E40").Value = "  +3.46%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "6.443"
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("D42").Value = "2.053"
$ws.Range("E42").Value = "  -2.57%  "
$ws.Range("D43").Value = "0.8834"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").Value = "0.4422"
$ws.Range("E44").Value = "  -0.83%  "
$ws.Range("D45").Value = "107.69"
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("D46").Value = "7.876"
$ws.Range("E46").Value = "  -3.69%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "980.79"
$ws.Range("E48").Value = "  -2.11%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "36.12"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.1237"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("D51").Value = "9.173"
$ws.Range("E51").Value = "  -2.05%  "
